# Updates the cryptocurrency price/volume snapshot in Sheet1 (cols B-E, rows 2-51)
# to match the latest scrape. Price cells in column D are written with a
# leading apostrophe (quote-prefix) so numeric-looking strings such as
# "0.999" or "5.21" are stored as text -- matching the original inlineStr
# text cells -- rather than being auto-coerced to Excel numbers. The style
# is reset to "Normal" right after so the cell keeps its original (default)
# style index instead of picking up a new "@" text-format style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.225.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.17%  '

$ws.Range("D3").Value = "'2.445.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'582.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.38%  '

$ws.Range("D6").Value = "'142.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = "'0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("D9").Value = "'2.439.89"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.41%  '

$ws.Range("E10").Value = '  +1.26%  '

$ws.Range("E11").Value = '  +2.61%  '

$ws.Range("D12").Value = "'5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").Value = "'0.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.81%  '

$ws.Range("D14").Value = "'26.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").Value = "'0.0000177"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.51%  '

$ws.Range("D16").Value = "'2.895.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.61%  '

$ws.Range("D17").Value = "'62.197.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").Value = "'2.442.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").Value = "'10.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.51%  '

$ws.Range("D20").Value = "'7.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").Value = "'326.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.88%  '

$ws.Range("D22").Value = "'4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '

$ws.Range("D23").Value = "'1.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.98%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = "'65.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("D26").Value = "'9.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.34%  '

$ws.Range("D27").Value = "'590.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.91%  '

$ws.Range("D28").Value = "'0.0₃0966"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.54%  '

$ws.Range("D29").Value = "'2.569.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.27%  '

$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").Value = "'1.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.94%  '

$ws.Range("D32").Value = "'7.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.80%  '

$ws.Range("E33").Value = '  +2.27%  '

$ws.Range("E34").Value = '  +0.91%  '

$ws.Range("D35").Value = "'4.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.88%  '

$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").Value = "'1.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.28%  '

$ws.Range("D38").Value = "'0.377"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.95%  '

$ws.Range("D39").Value = "'153.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.45%  '

$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D40").Value = "'18.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.38%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = "'5.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.81%  '

$ws.Range("D42").Value = "'42.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.87%  '

$ws.Range("D43").Value = "'1.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.43%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").Value = "'2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.28%  '

$ws.Range("D46").Value = "'142.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.77%  '

$ws.Range("D47").Value = "'3.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.19%  '

$ws.Range("D48").Value = "'0.0₆0253"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +14.24%  '

$ws.Range("D49").Value = "'0.605"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.99%  '

$ws.Range("D50").Value = "'0.0520"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '

$ws.Range("D51").Value = "'19.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.36%  '

